$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -15.67190036890694
$ws.Range("C2").Value = 1.322685100424217
$ws.Range("D2").Value = -15.67190036890694
$ws.Range("E2").Value = -15.67190036890694
$ws.Range("F2").Value = -15.67190036890694
$ws.Range("G2").Value = -15.67190036890694
$ws.Range("H2").Value = -15.67190036890694
$ws.Range("I2").Value = -15.67190036890694
$ws.Range("J2").Value = -15.67190036890694
$ws.Range("K2").Value = -15.67190036890694
$ws.Range("B3").Value = -15.67190036890694
$ws.Range("C3").Value = -15.67190036890694
$ws.Range("D3").Value = -15.67190036890694
$ws.Range("E3").Value = -15.67190036890694
$ws.Range("F3").Value = -15.67190036890694
$ws.Range("G3").Value = -15.67190036890694
$ws.Range("H3").Value = -15.67190036890694
$ws.Range("I3").Value = 1.222266083512366
$ws.Range("J3").Value = -15.67190036890694
$ws.Range("K3").Value = -15.67190036890694
$ws.Range("B4").Value = -15.67190036890694
$ws.Range("C4").Value = 1.655485611816391
$ws.Range("D4").Value = 4.321901841218952
$ws.Range("E4").Value = -15.67190036890694
$ws.Range("F4").Value = 3.52937404861395
$ws.Range("G4").Value = -15.67190036890694
$ws.Range("H4").Value = 1.891124555219899
$ws.Range("I4").Value = -15.67190036890694
$ws.Range("J4").Value = 2.491302129701983
$ws.Range("K4").Value = -15.67190036890694
$ws.Range("B5").Value = -15.67190036890694
$ws.Range("C5").Value = 1.824493132820061
$ws.Range("D5").Value = -15.67190036890694
$ws.Range("E5").Value = -15.67190036890694
$ws.Range("F5").Value = -15.67190036890694
$ws.Range("G5").Value = 3.442159650634955
$ws.Range("H5").Value = -15.67190036890694
$ws.Range("I5").Value = -15.67190036890694
$ws.Range("J5").Value = -15.67190036890694
$ws.Range("K5").Value = -15.67190036890694
$ws.Range("B6").Value = -15.67190036890694
$ws.Range("C6").Value = -15.67190036890694
$ws.Range("D6").Value = -15.67190036890694
$ws.Range("E6").Value = -15.67190036890694
$ws.Range("F6").Value = -15.67190036890694
$ws.Range("G6").Value = -15.67190036890694
$ws.Range("H6").Value = -15.67190036890694
$ws.Range("I6").Value = -15.67190036890694
$ws.Range("J6").Value = -15.67190036890694
$ws.Range("K6").Value = -15.67190036890694
$ws.Range("B7").Value = 3.01471300119778
$ws.Range("C7").Value = -15.67190036890694
$ws.Range("D7").Value = -15.67190036890694
$ws.Range("E7").Value = -15.67190036890694
$ws.Range("F7").Value = -15.67190036890694
$ws.Range("G7").Value = -15.67190036890694
$ws.Range("H7").Value = -15.67190036890694
$ws.Range("I7").Value = -15.67190036890694
$ws.Range("J7").Value = -15.67190036890694
$ws.Range("K7").Value = -15.67190036890694
$ws.Range("B8").Value = -15.67190036890694
$ws.Range("C8").Value = -15.67190036890694
$ws.Range("D8").Value = -15.67190036890694
$ws.Range("E8").Value = 2.421316004727642
$ws.Range("F8").Value = -15.67190036890694
$ws.Range("G8").Value = -15.67190036890694
$ws.Range("H8").Value = -15.67190036890694
$ws.Range("I8").Value = -15.67190036890694
$ws.Range("J8").Value = -15.67190036890694
$ws.Range("K8").Value = -15.67190036890694
$ws.Range("B9").Value = 3.575028068274172
$ws.Range("C9").Value = -15.67190036890694
$ws.Range("D9").Value = -15.67190036890694
$ws.Range("E9").Value = -15.67190036890694
$ws.Range("F9").Value = -15.67190036890694
$ws.Range("G9").Value = -15.67190036890694
$ws.Range("H9").Value = -15.67190036890694
$ws.Range("I9").Value = -15.67190036890694
$ws.Range("J9").Value = -15.67190036890694
$ws.Range("K9").Value = -15.67190036890694
$ws.Range("B10").Value = -15.67190036890694
$ws.Range("C10").Value = -15.67190036890694
$ws.Range("D10").Value = -15.67190036890694
$ws.Range("E10").Value = -15.67190036890694
$ws.Range("F10").Value = -15.67190036890694
$ws.Range("G10").Value = -15.67190036890694
$ws.Range("H10").Value = -15.67190036890694
$ws.Range("I10").Value = 0.6262517199027761
$ws.Range("J10").Value = -15.67190036890694
$ws.Range("K10").Value = 1.961215483099366
$ws.Range("B11").Value = -15.67190036890694
$ws.Range("C11").Value = -15.67190036890694
$ws.Range("D11").Value = -15.67190036890694
$ws.Range("E11").Value = 2.250759520608276
$ws.Range("F11").Value = -15.67190036890694
$ws.Range("G11").Value = 1.814471168983895
$ws.Range("H11").Value = -15.67190036890694
$ws.Range("I11").Value = -15.67190036890694
$ws.Range("J11").Value = -15.67190036890694
$ws.Range("K11").Value = 1.416359927672465
$ws.Range("B12").Value = -15.67190036890694
$ws.Range("C12").Value = -15.67190036890694
$ws.Range("D12").Value = -15.67190036890694
$ws.Range("E12").Value = -15.67190036890694
$ws.Range("F12").Value = -15.67190036890694
$ws.Range("G12").Value = -15.67190036890694
$ws.Range("H12").Value = -15.67190036890694
$ws.Range("I12").Value = -15.67190036890694
$ws.Range("J12").Value = -15.67190036890694
$ws.Range("K12").Value = -15.67190036890694
$ws.Range("B13").Value = -15.67190036890694
$ws.Range("C13").Value = -15.67190036890694
$ws.Range("D13").Value = -15.67190036890694
$ws.Range("E13").Value = 2.280300539856549
$ws.Range("F13").Value = -15.67190036890694
$ws.Range("G13").Value = -15.67190036890694
$ws.Range("H13").Value = -15.67190036890694
$ws.Range("I13").Value = -15.67190036890694
$ws.Range("J13").Value = 1.187258262193424
$ws.Range("K13").Value = 2.349280243737554
$ws.Range("B14").Value = -15.67190036890694
$ws.Range("C14").Value = -15.67190036890694
$ws.Range("D14").Value = -15.67190036890694
$ws.Range("E14").Value = -15.67190036890694
$ws.Range("F14").Value = -15.67190036890694
$ws.Range("G14").Value = -15.67190036890694
$ws.Range("H14").Value = -15.67190036890694
$ws.Range("I14").Value = -15.67190036890694
$ws.Range("J14").Value = -15.67190036890694
$ws.Range("K14").Value = 2.023851272306104
$ws.Range("B15").Value = -15.67190036890694
$ws.Range("C15").Value = -15.67190036890694
$ws.Range("D15").Value = -15.67190036890694
$ws.Range("E15").Value = -15.67190036890694
$ws.Range("F15").Value = -15.67190036890694
$ws.Range("G15").Value = -15.67190036890694
$ws.Range("H15").Value = -15.67190036890694
$ws.Range("I15").Value = -15.67190036890694
$ws.Range("J15").Value = -15.67190036890694
$ws.Range("K15").Value = -15.67190036890694
$ws.Range("B16").Value = -15.67190036890694
$ws.Range("C16").Value = -15.67190036890694
$ws.Range("D16").Value = -15.67190036890694
$ws.Range("E16").Value = -15.67190036890694
$ws.Range("F16").Value = -15.67190036890694
$ws.Range("G16").Value = -15.67190036890694
$ws.Range("H16").Value = -15.67190036890694
$ws.Range("I16").Value = -15.67190036890694
$ws.Range("J16").Value = 2.040203401124728
$ws.Range("K16").Value = -15.67190036890694
$ws.Range("B17").Value = -15.67190036890694
$ws.Range("C17").Value = 2.112025069209759
$ws.Range("D17").Value = -15.67190036890694
$ws.Range("E17").Value = -15.67190036890694
$ws.Range("F17").Value = -15.67190036890694
$ws.Range("G17").Value = -15.67190036890694
$ws.Range("H17").Value = 2.10098763243221
$ws.Range("I17").Value = 1.287116707216543
$ws.Range("J17").Value = 2.198494454675969
$ws.Range("K17").Value = -15.67190036890694
$ws.Range("B18").Value = -15.67190036890694
$ws.Range("C18").Value = -15.67190036890694
$ws.Range("D18").Value = -15.67190036890694
$ws.Range("E18").Value = -15.67190036890694
$ws.Range("F18").Value = -15.67190036890694
$ws.Range("G18").Value = -15.67190036890694
$ws.Range("H18").Value = 1.894398994588796
$ws.Range("I18").Value = 0.5807835405643457
$ws.Range("J18").Value = 1.764103509232836
$ws.Range("K18").Value = -15.67190036890694
$ws.Range("B19").Value = -15.67190036890694
$ws.Range("C19").Value = -15.67190036890694
$ws.Range("D19").Value = -15.67190036890694
$ws.Range("E19").Value = -15.67190036890694
$ws.Range("F19").Value = -15.67190036890694
$ws.Range("G19").Value = -15.67190036890694
$ws.Range("H19").Value = 1.675175941634138
$ws.Range("I19").Value = 1.655485611816391
$ws.Range("J19").Value = -15.67190036890694
$ws.Range("K19").Value = -15.67190036890694
$ws.Range("B20").Value = -15.67190036890694
$ws.Range("C20").Value = 1.259910552827073
$ws.Range("D20").Value = -15.67190036890694
$ws.Range("E20").Value = -15.67190036890694
$ws.Range("F20").Value = 3.079503122033454
$ws.Range("G20").Value = -15.67190036890694
$ws.Range("H20").Value = 1.241327994492405
$ws.Range("I20").Value = 3.175796262911361
$ws.Range("J20").Value = -15.67190036890694
$ws.Range("K20").Value = 2.095695724724557
$ws.Range("B21").Value = -15.67190036890694
$ws.Range("C21").Value = 2.031577436402553
$ws.Range("D21").Value = -15.67190036890694
$ws.Range("E21").Value = 2.329403358267016
$ws.Range("F21").Value = -15.67190036890694
$ws.Range("G21").Value = 2.488842634882334
$ws.Range("H21").Value = 1.446280315403681
$ws.Range("I21").Value = -15.67190036890694
$ws.Range("J21").Value = -15.67190036890694
$ws.Range("K21").Value = -15.67190036890694
